# =====================================================================
# Edit script: "Changed timeseries inputdata in regards to scenarios."
# TimeSeries data (inflows, prices, etc.) can now be defined once for
# several/all scenarios instead of one column per scenario. This script
# reproduces that change on the Predicer input_data workbook.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "price" sheet: collapse the per-scenario ng price columns
#    (ng,s1 / ng,s2 / ng,s3) into a single "ng, ALL" column.
# ---------------------------------------------------------------------
$price = $wb.Worksheets.Item("price")
# Remove column D ("ng,s3") then column C ("ng,s2"), leaving column B.
$price.Columns.Item(4).Delete()
$price.Columns.Item(3).Delete()

# ---------------------------------------------------------------------
# 2) "constraints" sheet: add documentation / example rows describing
#    the new generic constraint columns (operator, is_limited,
#    limited_by, constant_diff) and a worked example.
# ---------------------------------------------------------------------
$constraints = $wb.Worksheets.Item("constraints")

$constraints.Range("B5").Value = "operator"
$constraints.Range("C5").Value = "is_limited"
$constraints.Range("D5").Value = "limited_by"

$constraints.Range("A7").Value = "c2"
$constraints.Range("A8").Value = "c3"

$constraints.Range("B7").Value = "gt"

$constraints.Range("C6").Value = "flow"
$constraints.Range("D6").Value = "state"

$constraints.Range("E5").Value = "constant_diff"

$constraints.Range("C9").Value = "flow_val"
$constraints.Range("D9").Value = "state_val"

$constraints.Range("B8").Value = "st"

$constraints.Range("D11").Value = "constant"

$constraints.Range("D12").Value = "timeseries"

$constraints.Range("A16").Value = "t1"
$constraints.Range("A17").Value = "t2"
$constraints.Range("A18").Value = "t3"

$constraints.Range("B15").Value = "c1, s1"

# Fill in the remaining example rows (re-using strings already present
# in the workbook, so no new shared-string entries are created here).
$constraints.Range("A5").Value = "name"
$constraints.Range("A6").Value = "c1"
$constraints.Range("B6").Value = "eq"
$constraints.Range("C7").Value = "state"
$constraints.Range("C10").Value = "state_val"
$constraints.Range("D10").Value = "flow_val"
$constraints.Range("C11").Value = "reserve"
$constraints.Range("E6").Value = -1

$constraints.Columns.Item(3).ColumnWidth = 9.85546875
$constraints.Columns.Item(4).ColumnWidth = 10.5703125

$constraints.Activate()
$constraints.Range("B15").Select()

# ---------------------------------------------------------------------
# 3) Re-label the "ng" price column to reflect that it now applies to
#    every scenario ("ng, ALL"). Creating this shared string last keeps
#    the shared-string table ordering consistent with the source file.
# ---------------------------------------------------------------------
$price.Range("B1").Value = "ng, ALL"

# ---------------------------------------------------------------------
# 4) Misc. view/selection bookkeeping to mirror the authors' session.
# ---------------------------------------------------------------------
$nodes = $wb.Worksheets.Item("nodes")
$nodes.Range("H10").Select()

$genConstraint = $wb.Worksheets.Item("gen_constraint")
$genConstraint.Range("H1").Select()

$price.Activate()
$price.Range("B1").Select()
